$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 135, pushing the existing
# rows 135-143 down to 136-144 (dimension grows from R143 to R144).
$ws.Rows.Item(135).Insert()

# Populate the newly inserted row 135 with the new record's data.
$ws.Cells.Item(135, 1).Value = 7
$ws.Cells.Item(135, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(135, 3).Value = "Ñuble"
$ws.Cells.Item(135, 4).Value = 44516
$ws.Cells.Item(135, 5).Value = 16
$ws.Cells.Item(135, 6).Value = 100112032
$ws.Cells.Item(135, 7).Value = "Zapallo italiano"
$ws.Cells.Item(135, 8).Value = "Sin especificar"
$ws.Cells.Item(135, 9).Value = "Primera"
$ws.Cells.Item(135, 10).Value = 100
$ws.Cells.Item(135, 11).Value = 8000
$ws.Cells.Item(135, 12).Value = 9000
$ws.Cells.Item(135, 13).Value = 8500
$ws.Cells.Item(135, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(135, 15).Value = "Región del Maule"
$ws.Cells.Item(135, 16).Value = 142
$ws.Cells.Item(135, 17).Value = 60
$ws.Cells.Item(135, 18).Value = "Hortaliza"
